# "Add files via upload" — apply the recorded edits to "Diseño del enlace.xlsx".
# Sheet2 (the active sheet) gains a new Sx/BW column plus an antenna-gain
# comparison block in H:J, F2:F4 switch from a placeholder label to actual
# dBm readings, E3/E4 become live hyperlinks, and the window view/zoom saved
# state is refreshed to match the author's last on-screen selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- F2:F4 — replace the "a" placeholder text with the real computed values ---
$ws.Range("F2").Value = -75
$ws.Range("F3").Value = -77.239999999999995
$ws.Range("F4").Value = -77

# --- New header cells on row 1 ---
$ws.Range("F1").Value = "Sx @ BW 200MHz"
$ws.Range("H1").Value = "RX_Sens(dBm) = –174(dBm) + 10*LOG[BW(Hz)] + NF(dB) + S/N(dB)"

# --- E3 / E4 become real hyperlinks (text already held the target URL) ---
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.olifantasia.com/gnuradio/usrp/files/datasheets/ds_n200series.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://www.ettus.com/content/files/07495_Ettus_N200-210_DS_Flyer_HR_1.pdf") | Out-Null
# match the existing hyperlink-style cells (E2, D19, D20)
$ws.Range("E3").Style = $ws.Range("E2").Style
$ws.Range("E4").Style = $ws.Range("E2").Style

# --- New comparison block, columns H:J, rows 16-22 ---
$ws.Range("H16").Value = "HackRF One "
$ws.Range("I16").Value = "CBX"
$ws.Range("J16").Value = "UBX"

$ws.Range("H17").Value = 1.6
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 13

$ws.Range("H18").Value = 35
$ws.Range("I18").Value = 410
$ws.Range("J18").Value = 270

$ws.Range("H19").Value = "-"
$ws.Range("I19").Value = 1.6
$ws.Range("J19").Value = "-"

$ws.Range("H20").Value = "-"
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 3

$ws.Range("H21").Value = 35
$ws.Range("I21").Value = 410
$ws.Range("J21").Value = 270

$ws.Range("H22").Value = "-"
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = 5.2

# --- Column F needs to be wide enough to show "Sx @ BW 200MHz" ---
$ws.Columns.Item(6).AutoFit() | Out-Null

# --- Window/view state: scrolled one column right, E18 selected ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E18").Select() | Out-Null
